# "Generate Report for Handoff"
#
# Refreshes the localization-status report: the handoff run that produced
# this workbook generated a new HO xliff for bafcc2e9-6684-44b8-87dc-675f32e05915.md
# and picked up new "Latest Handoff Datetime" stamps for the zh-cn/de-de
# targets. Stamp the new timestamps onto the affected cells.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for bafcc2e9-... (row 7)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-17 04:38:08"

# zh-cn sheet: "Latest Handoff Datetime" for bafcc2e9-... (row 7)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-17 04:37:59"

# de-de sheet: "Latest Handoff Datetime" for 8ba283ea-... (row 6) and
# bafcc2e9-... (row 7)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-17 04:38:08"
$wsDeDe.Range("H7").Value = "2016-08-17 04:38:08"
